$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates derived from the refreshed crypto price feed.
# Numeric-looking price strings are prefixed with a literal apostrophe so
# Excel stores them as text (matching the original inlineStr cells) instead
# of silently converting them to the Number type.

$ws.Range('D2').Value = '27.902.06'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '1.632.13'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''211.56'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '''23.52'
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').Value = '1.627.18'
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = '''0.564'
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '27.903.40'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '''229.19'
$ws.Range('E18').Value = '  -1.59%  '
$ws.Range('D19').Value = '''7.68'
$ws.Range('E19').Value = '  +1.81%  '
$ws.Range('D20').Value = '0.0₃0719'
$ws.Range('E20').Value = '  -0.29%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('E22').Value = '  -0.88%  '
$ws.Range('E23').Value = '  -3.41%  '
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('D25').Value = '''154.59'
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('E26').Value = '  -0.73%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').Value = '''15.53'
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').Value = '''3.42'
$ws.Range('E32').Value = '  +0.94%  '
$ws.Range('E33').Value = '  +0.76%  '
$ws.Range('D34').Value = '1.393.67'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('E36').Value = '  +10.14%  '
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('E38').Value = '  +1.18%  '
$ws.Range('D39').Value = '''0.558'
$ws.Range('E39').Value = '  -0.84%  '
$ws.Range('D40').Value = '''0.848'
$ws.Range('E40').Value = '  -3.25%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('D44').Value = '''65.81'
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('D45').Value = '''5.43'
$ws.Range('E45').Value = '  -2.00%  '
$ws.Range('D47').Value = '''2.14'
$ws.Range('E47').Value = '  -3.29%  '
$ws.Range('D48').Value = '''88.64'
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('D49').Value = '''0.102'
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0101'
$ws.Range('E50').Value = '  -3.56%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.0505'
$ws.Range('E51').Value = '  -0.11%  '
